# Auto-generated Excel COM-interop script to update Goblin_Profits market data
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets) per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1275.875
$ws.Range("J129").Value = 2308
$ws.Range("L129").Value = 6924
$ws.Range("N129").Value = -16924
$ws.Range("H135").Value = 602.9677
$ws.Range("I135").Value = 589.73334
$ws.Range("K135").Value = 5307.60006
$ws.Range("M135").Value = -2772.60006

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 59668.555
$ws.Range("I32").Value = 62025.254
$ws.Range("K32").Value = 62025.254
$ws.Range("M32").Value = -61738.254
$ws.Range("H63").Value = 7041.4707
$ws.Range("I63").Value = 4080.5
$ws.Range("K63").Value = 4080.5
$ws.Range("M63").Value = -3394.5
$ws.Range("H66").Value = 7041.4707
$ws.Range("I66").Value = 4080.5
$ws.Range("K66").Value = 20402.5
$ws.Range("M66").Value = -16970.5
$ws.Range("H102").Value = 7139.8667
$ws.Range("J102").Value = 7333.3335
$ws.Range("L102").Value = 7333.3335
$ws.Range("N102").Value = -10577.3335
$ws.Range("H110").Value = 3101.5
$ws.Range("I110").Value = 830.3333
$ws.Range("K110").Value = 830.3333
$ws.Range("M110").Value = 1214.6667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7851.1875
$ws.Range("I107").Value = 7386.5
$ws.Range("K107").Value = 7386.5
$ws.Range("M107").Value = -5466.5
$ws.Range("H135").Value = 105748.5
$ws.Range("J135").Value = 107098.4
$ws.Range("L135").Value = 107098.4
$ws.Range("N135").Value = -117238.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 11895
$ws.Range("I57").Value = 11895
$ws.Range("K57").Value = 11895
$ws.Range("M57").Value = -11335
$ws.Range("H94").Value = 1672.2307
$ws.Range("I94").Value = 1054.3334
$ws.Range("J94").Value = 1857.6
$ws.Range("K94").Value = 1054.3334
$ws.Range("L94").Value = 1857.6
$ws.Range("M94").Value = -603.3334
$ws.Range("N94").Value = -2759.6
$ws.Range("H132").Value = 4554.074
$ws.Range("I132").Value = 1639.0588
$ws.Range("K132").Value = 4917.1764
$ws.Range("M132").Value = -2387.1764
$ws.Range("H135").Value = 102944
$ws.Range("J135").Value = 102944
$ws.Range("L135").Value = 102944
$ws.Range("N135").Value = -113084

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 98.083336
$ws.Range("J2").Value = 263
$ws.Range("L2").Value = 1578
$ws.Range("N2").Value = -1804
$ws.Range("H6").Value = 27.7
$ws.Range("I6").Value = 33.5
$ws.Range("J6").Value = 4.5
$ws.Range("K6").Value = 100.5
$ws.Range("L6").Value = 13.5
$ws.Range("M6").Value = 12.5
$ws.Range("N6").Value = -239.5
$ws.Range("H7").Value = 44301.707
$ws.Range("I7").Value = 318.66666
$ws.Range("J7").Value = 68292.45
$ws.Range("K7").Value = 955.9999799999999
$ws.Range("L7").Value = 204877.35
$ws.Range("M7").Value = -843.9999799999999
$ws.Range("N7").Value = -205101.35
$ws.Range("H18").Value = 532.8
$ws.Range("I18").Value = 532.8
$ws.Range("K18").Value = 1598.4
$ws.Range("M18").Value = -1429.4
$ws.Range("H31").Value = 450
$ws.Range("I31").Value = 450
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1350
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1062
$ws.Range("N31").ClearContents()
$ws.Range("H68").Value = 1902.2222
$ws.Range("I68").Value = 1198.75
$ws.Range("J68").Value = 2465
$ws.Range("K68").Value = 3596.25
$ws.Range("L68").Value = 7395
$ws.Range("M68").Value = -2785.25
$ws.Range("N68").Value = -9017
$ws.Range("H71").Value = 1902.2222
$ws.Range("I71").Value = 1198.75
$ws.Range("J71").Value = 2465
$ws.Range("K71").Value = 10788.75
$ws.Range("L71").Value = 22185
$ws.Range("M71").Value = -6732.75
$ws.Range("N71").Value = -30297
$ws.Range("H95").Value = 18999.889
$ws.Range("J95").Value = 18999.889
$ws.Range("L95").Value = 56999.667
$ws.Range("N95").Value = -61117.667
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H114").Value = 2049.375
$ws.Range("I114").Value = 1484.1666
$ws.Range("K114").Value = 4452.4998
$ws.Range("M114").Value = -1198.4998
$ws.Range("H129").Value = 2714.25
$ws.Range("I129").Value = 888.2727
$ws.Range("K129").Value = 2664.8181
$ws.Range("M129").Value = 2335.1819
$ws.Range("H134").Value = 8810.210999999999
$ws.Range("I134").Value = 2466
$ws.Range("J134").Value = 9999.75
$ws.Range("K134").Value = 7398
$ws.Range("L134").Value = 29999.25
$ws.Range("M134").Value = -2328
$ws.Range("N134").Value = -40139.25
$ws.Range("H139").Value = 3968.7354
$ws.Range("J139").Value = 4117.5293
$ws.Range("L139").Value = 12352.5879
$ws.Range("N139").Value = -22632.5879
$ws.Range("H141").Value = 142861550
$ws.Range("I141").Value = 250003380
$ws.Range("K141").Value = 750010140
$ws.Range("M141").Value = -750004960

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1130.138
$ws.Range("I102").Value = 912.9583
$ws.Range("K102").Value = 912.9583
$ws.Range("M102").Value = 709.0417
$ws.Range("H107").Value = 1995.5333
$ws.Range("I107").Value = 247.57143
$ws.Range("J107").Value = 3525
$ws.Range("K107").Value = 247.57143
$ws.Range("L107").Value = 3525
$ws.Range("M107").Value = 1672.42857
$ws.Range("N107").Value = -7365
$ws.Range("H132").Value = 100002860
$ws.Range("I132").Value = 142859090
$ws.Range("K132").Value = 428577270
$ws.Range("M132").Value = -428574740

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2400
$ws.Range("I3").Value = 1650
$ws.Range("J3").Value = 3150
$ws.Range("K3").Value = 1650
$ws.Range("L3").Value = 3150
$ws.Range("M3").Value = -1538
$ws.Range("N3").Value = -3374
$ws.Range("H14").Value = 3850
$ws.Range("J14").Value = 2200
$ws.Range("L14").Value = 2200
$ws.Range("N14").Value = -2544
$ws.Range("H15").Value = 2400
$ws.Range("I15").Value = 1650
$ws.Range("J15").Value = 3150
$ws.Range("K15").Value = 1650
$ws.Range("L15").Value = 3150
$ws.Range("M15").Value = -1480
$ws.Range("N15").Value = -3490
$ws.Range("H20").Value = 292127
$ws.Range("I20").Value = 669666.7
$ws.Range("J20").Value = 8972.25
$ws.Range("K20").Value = 669666.7
$ws.Range("L20").Value = 8972.25
$ws.Range("M20").Value = -669440.7
$ws.Range("N20").Value = -9424.25
$ws.Range("H21").Value = 14900
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 14900
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 14900
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -15248
$ws.Range("H22").Value = 3231.1765
$ws.Range("I22").Value = 2658
$ws.Range("J22").Value = 3543.818
$ws.Range("K22").Value = 2658
$ws.Range("L22").Value = 3543.818
$ws.Range("M22").Value = -2363
$ws.Range("N22").Value = -4133.818
$ws.Range("H27").Value = 3231.1765
$ws.Range("I27").Value = 2658
$ws.Range("J27").Value = 3543.818
$ws.Range("K27").Value = 2658
$ws.Range("L27").Value = 3543.818
$ws.Range("M27").Value = -2551
$ws.Range("N27").Value = -3757.818
$ws.Range("H40").Value = 4056.4062
$ws.Range("I40").Value = 2991.182
$ws.Range("K40").Value = 2991.182
$ws.Range("M40").Value = -2855.182
$ws.Range("H46").Value = 1037.2222
$ws.Range("I46").Value = 435.2857
$ws.Range("K46").Value = 435.2857
$ws.Range("M46").Value = -247.2857
$ws.Range("H93").Value = 4668.3057
$ws.Range("I93").Value = 2791.8823
$ws.Range("J93").Value = 6347.2104
$ws.Range("K93").Value = 2791.8823
$ws.Range("L93").Value = 6347.2104
$ws.Range("M93").Value = -1543.8823
$ws.Range("N93").Value = -8843.2104
$ws.Range("H122").Value = 4522.75
$ws.Range("I122").Value = 3961.182
$ws.Range("J122").Value = 5758.2
$ws.Range("K122").Value = 11883.546
$ws.Range("L122").Value = 17274.6
$ws.Range("M122").Value = -9433.545999999998
$ws.Range("N122").Value = -22174.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11114403
$ws.Range("I132").Value = 15875547
$ws.Range("K132").Value = 47626641
$ws.Range("M132").Value = -47624111

